# Auto-generated Excel COM-interop script applying cell value updates
# derived from the OOXML diff for Marilith_Profits workbook.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 2 (diff @@ -727)
$ws.Range("H2").Value = 4057
$ws.Range("J2").Value = 5179.8
$ws.Range("L2").Value = 5179.8
$ws.Range("N2").Value = -5405.8
# row 21 (diff @@ -1691)
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# row 23 (diff @@ -1789)
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# row 29 (diff @@ -2080)
$ws.Range("H29").Value = 697.1667
$ws.Range("I29").Value = 136.6
$ws.Range("K29").Value = 409.8
$ws.Range("M29").Value = -128.8
# row 38 (diff @@ -2518)
$ws.Range("H38").Value = 570.36365
$ws.Range("I38").Value = 221.75
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 665.25
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = -293.25
$ws.Range("N38").Value = -5244
# row 58 (diff @@ -3492)
$ws.Range("H58").Value = 1577.75
$ws.Range("J58").Value = 2126.8333
$ws.Range("L58").Value = 6380.499899999999
$ws.Range("N58").Value = -6680.499899999999
# row 62 (diff @@ -3685)
$ws.Range("H62").Value = 2538.6
# row 65 (diff @@ -3835)
$ws.Range("H65").Value = 2538.6
# row 87 (diff @@ -4931)
$ws.Range("H87").Value = 49999.5
$ws.Range("J87").Value = 49999.5
$ws.Range("L87").Value = 49999.5
$ws.Range("N87").Value = -52495.5
# row 90 (diff @@ -5081)
$ws.Range("H90").Value = 49999.5
$ws.Range("J90").Value = 49999.5
$ws.Range("L90").Value = 149998.5
$ws.Range("N90").Value = -162478.5
# row 129 (diff @@ -6980)
$ws.Range("H129").Value = 1777.85
$ws.Range("I129").Value = 958.0714
$ws.Range("J129").Value = 3690.6667
$ws.Range("K129").Value = 2874.2142
$ws.Range("L129").Value = 11072.0001
$ws.Range("M129").Value = 2125.7858
$ws.Range("N129").Value = -21072.0001
# row 135 (diff @@ -7274)
$ws.Range("H135").Value = 797.3333
$ws.Range("J135").Value = 1389.5
$ws.Range("L135").Value = 12505.5
$ws.Range("N135").Value = -17575.5
# row 137 (diff @@ -7372)
$ws.Range("H137").Value = 2395
$ws.Range("J137").Value = 2597.5
$ws.Range("L137").Value = 7792.5
$ws.Range("N137").Value = -12892.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 4 (diff @@ -7815)
$ws.Range("H4").Value = 834.2222
$ws.Range("I4").Value = 844.75
$ws.Range("K4").Value = 844.75
$ws.Range("M4").Value = -728.75
# row 6 (diff @@ -7919)
$ws.Range("H6").Value = 999
$ws.Range("I6").Value = 999
$ws.Range("K6").Value = 999
$ws.Range("M6").Value = -826
# row 61 (diff @@ -10614)
$ws.Range("H61").Value = 3470.6667
$ws.Range("I61").Value = 3164.8
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3164.8
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2952.8
$ws.Range("N61").Value = -5424
# row 63 (diff @@ -10712)
$ws.Range("H63").Value = 2307.3333
$ws.Range("I63").Value = 1968.8
$ws.Range("K63").Value = 1968.8
$ws.Range("M63").Value = -1282.8
# row 66 (diff @@ -10859)
$ws.Range("H66").Value = 2307.3333
$ws.Range("I66").Value = 1968.8
$ws.Range("K66").Value = 9844
$ws.Range("M66").Value = -6412
# row 136 (diff @@ -14217)
$ws.Range("H136").Value = 3470.6667
$ws.Range("I136").Value = 3164.8
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9494.400000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6944.400000000001
$ws.Range("N136").Value = -20100

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row 86 (diff @@ -18661)
$ws.Range("H86").Value = 1268
$ws.Range("I86").Value = 2727.6667
$ws.Range("J86").Value = 642.4286
$ws.Range("K86").Value = 2727.6667
$ws.Range("L86").Value = 642.4286
$ws.Range("M86").Value = -1604.6667
$ws.Range("N86").Value = -2888.4286
# row 89 (diff @@ -18808)
$ws.Range("H89").Value = 1268
$ws.Range("I89").Value = 2727.6667
$ws.Range("J89").Value = 642.4286
$ws.Range("K89").Value = 13638.3335
$ws.Range("L89").Value = 3212.143
$ws.Range("M89").Value = -8022.333500000001
$ws.Range("N89").Value = -14444.143
# row 94 (diff @@ -19050)
$ws.Range("H94").Value = 2739.4666
$ws.Range("I94").Value = 2853.2307
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 2853.2307
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -2402.2307
$ws.Range("N94").Value = -2902

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 12 (diff @@ -21950)
$ws.Range("H12").Value = 19711.572
$ws.Range("J12").Value = 22830.166
$ws.Range("L12").Value = 22830.166
$ws.Range("N12").Value = -23170.166
# row 86 (diff @@ -25609)
$ws.Range("H86").Value = 6000
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
# row 89 (diff @@ -25756)
$ws.Range("H89").Value = 6000
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 4 (diff @@ -28515)
$ws.Range("H4").Value = 1891.875
$ws.Range("I4").Value = 2052.8572
$ws.Range("J4").Value = 1666.5
$ws.Range("K4").Value = 6158.571599999999
$ws.Range("L4").Value = 4999.5
$ws.Range("M4").Value = -6046.571599999999
$ws.Range("N4").Value = -5223.5
# row 37 (diff @@ -30180)
$ws.Range("H37").Value = 96666.664
$ws.Range("J37").Value = 96666.664
$ws.Range("L37").Value = 289999.992
$ws.Range("N37").Value = -290223.992
# row 137 (diff @@ -35134)
$ws.Range("H137").Value = 3569.6667
$ws.Range("I137").Value = 2906
$ws.Range("J137").Value = 4233.3335
$ws.Range("K137").Value = 8718
$ws.Range("L137").Value = 12700.0005
$ws.Range("M137").Value = -3618
$ws.Range("N137").Value = -22900.0005
# row 138 (diff @@ -35186)
$ws.Range("H138").Value = 4614.5835
$ws.Range("I138").Value = 4481.5
$ws.Range("K138").Value = 13444.5
$ws.Range("M138").Value = -8304.5
# row 139 (diff @@ -35238)
$ws.Range("H139").Value = 1156
$ws.Range("I139").Value = 1157.7693
$ws.Range("J139").Value = 1133
$ws.Range("K139").Value = 3473.3079
$ws.Range("L139").Value = 3399
$ws.Range("M139").Value = 1666.6921
$ws.Range("N139").Value = -13679

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 126 (diff @@ -41474)
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 40 (diff @@ -44187)
$ws.Range("H40").Value = 3333.6667
$ws.Range("I40").Value = 2748.5
$ws.Range("K40").Value = 2748.5
$ws.Range("M40").Value = -2612.5
# row 82 (diff @@ -46227)
$ws.Range("H82").Value = 2160
$ws.Range("I82").Value = 1700
$ws.Range("K82").Value = 1700
$ws.Range("M82").Value = -1339
# row 85 (diff @@ -46380)
$ws.Range("H85").Value = 2160
$ws.Range("I85").Value = 1700
$ws.Range("K85").Value = 1700
$ws.Range("M85").Value = -452

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 81 (diff @@ -53006)
$ws.Range("H81").Value = 299.5
$ws.Range("I81").Value = 299.5
$ws.Range("K81").Value = 599
$ws.Range("M81").Value = 462
# row 84 (diff @@ -53147)
$ws.Range("H84").Value = 299.5
$ws.Range("I84").Value = 299.5
$ws.Range("K84").Value = 2995
$ws.Range("M84").Value = 2309
# row 122 (diff @@ -54970)
$ws.Range("H122").Value = 2768
$ws.Range("I122").Value = 2496.75
$ws.Range("J122").Value = 3129.6667
$ws.Range("K122").Value = 7490.25
$ws.Range("L122").Value = 9389.000100000001
$ws.Range("M122").Value = -5040.25
$ws.Range("N122").Value = -14289.0001
# row 126 (diff @@ -55166)
$ws.Range("H126").Value = 1293.9166
$ws.Range("I126").Value = 1138.8182
$ws.Range("K126").Value = 3416.4546
$ws.Range("M126").Value = -946.4546
